$wb = $excel.ActiveWorkbook

# Add the new "Self Assessment Delegates" sheet after "Course Delegates"
$courseDelegatesSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $courseDelegatesSheet)
$ws.Name = "Self Assessment Delegates"

# Header row for the new sheet
$headers = @(
    "Self assessment name",
    "Last name",
    "First name",
    "Email",
    "PRN",
    "Role type",
    "Manager",
    "Base / office / place of work",
    "Base / office / place of work (Prompt 4)",
    "Contact telephone number",
    "Delegate ID",
    "Enrolled",
    "Last accessed",
    "Complete by",
    "Submitted",
    "Signed off",
    "Launches"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the header + blank data row into an Excel Table ("Table2")
$tableRange = $ws.Range("A1:Q2")
$listObj = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$listObj.Name = "Table2"

# Bold the "Self assessment name" column (first table column)
$listObj.ListColumns.Item(1).Range.Font.Bold = $true

# Match the original sheet's outline/grouping on the data row
$ws.Rows.Item(2).OutlineLevel = 1

# Make the new sheet the active tab, scrolled/selected like the source
$ws.Activate()
$ws.Range("K5").Select()

$excel.ActiveWindow.ScrollColumn = 8

Write-Output "Self Assessment Delegates sheet created"
